# Updated the results and edited the code so that it works regardless of
# the order of the quarters. This fills in the "Round_1" picks in column O
# (the mirrored / right-hand side of the bracket) to match the winners
# already chosen in column P for each corresponding match row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$picks = @{
    "O19" = "Emperor Penguin"
    "O21" = "Owl Monkey"
    "O23" = "Pacific Spiny Lumpsucker"
    "O25" = "Siamang"
    "O27" = "Bat-Eared Fox"
    "O29" = "Wolverine"
    "O31" = "Dyak Friut Bat"
    "O33" = "Greater Rhea"
}

foreach ($addr in $picks.Keys) {
    $ws.Range($addr).Value = $picks[$addr]
}

# Update the active selection to match where the editor was last working.
$ws.Range("O19").Select()
